$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-5 held the same three records but in the wrong order. Re-write each
# cell in rows 3-5 so that the records end up in the correct order:
#   new row 3 <- old row 5 data
#   new row 4 <- old row 3 data
#   new row 5 <- old row 4 data
# (row 4 only needs its A/B/E/F/G/H cells touched, because old rows 3 and 5
#  are the ones swapping places; the rest of row 4 was already correct.)

$ws.Range("A3").Value = 55773928
$ws.Range("B3").Value = 95717
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 220686
$ws.Range("F3").Value = 'Kambräken'
$ws.Range("G3").Value = 'Blechnum spicant'
$ws.Range("H3").Value = '(L.) Roth'
$ws.Range("K3").Value = ""
$ws.Range("P3").Value = 'Mellan Öster-Vike och Sjömyran, Ång'
$ws.Range("Q3").Value = 573419.838285814
$ws.Range("R3").Value = 7017474.295048638
$ws.Range("S3").Value = 10
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = '2015-08-30'
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = '2015-08-30'
$ws.Range("AH3").Value = 'Granskog'
$ws.Range("AN3").Value = 1
$ws.Range("AO3").Value = '1 substratenheter'
$ws.Range("AW3").Value = 'Magnus Johansson'
$ws.Range("AX3").Value = 'Magnus Johansson'
$ws.Range("AY3").Value = 'SCA Skog Naturvärdesinventering'
$ws.Range("A4").Value = 104427604
$ws.Range("B4").Value = 78570
$ws.Range("E4").Value = 2081
$ws.Range("F4").Value = 'Skrovellav'
$ws.Range("G4").Value = 'Lobaria scrobiculata'
$ws.Range("H4").Value = '(Scop.) DC.'
$ws.Range("A5").Value = 104427596
$ws.Range("B5").Value = 78569
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = 'Lunglav'
$ws.Range("G5").Value = 'Lobaria pulmonaria'
$ws.Range("H5").Value = '(L.) Hoffm.'
$ws.Range("K5").Value = ""
$ws.Range("P5").Value = 'Sollefteå, Ång'
$ws.Range("Q5").Value = 573169.9361146218
$ws.Range("R5").Value = 7017540.184384095
$ws.Range("S5").Value = 25
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = '2022-11-01'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = '2022-11-01'
$ws.Range("AH5").Value = ""
$ws.Range("AN5").Value = ""
$ws.Range("AO5").Value = ""
$ws.Range("AW5").Value = 'Erland Lindblad'
$ws.Range("AX5").Value = 'Erland Lindblad'
$ws.Range("AY5").Value = ""
